{"js": "// Update the invoice document:\n//  - document date 27.5.2024 -> 6.6.2024 (both occurrences)\n//  - item name KALTMANN MASCHINEN 2 -> YH9000AE\n//  - remove the \"14450.00\" price/sum figures (all three occurrences)\n//  - shorten the spelled-out sum text to just \"\u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\"\n\nconst body = context.document.body;\n\n// 1) Replace both occurrences of the old date with the new date.\nconst dateResults = body.search(\"27.5.2024\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"6.6.2024\", Word.InsertLocation.replace);\n}\n\n// 2) Replace the equipment name.\nconst nameResults = body.search(\"KALTMANN MASCHINEN 2\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < nameResults.items.length; i++) {\n  nameResults.items[i].insertText(\"YH9000AE\", Word.InsertLocation.replace);\n}\n\n// 3) Clear the \"14450.00\" amounts (unit price, line total, and the amount\n//    embedded in the \"\u0423\u0421\u042c\u041e\u0413\u041e\" summary cell).\nconst amountResults = body.search(\"14450.00\", { matchCase: true });\namountResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < amountResults.items.length; i++) {\n  amountResults.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\n\n// 4) Shorten the spelled-out sum text.\nconst sumResults = body.search(\n  \"\u0447\u043e\u0442\u0438\u0440\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0447\u043e\u0442\u0438\u0440\u0438\u0441\u0442\u0430 \u043f\u044f\u0442\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\",\n  { matchCase: true }\n);\nsumResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < sumResults.items.length; i++) {\n  sumResults.items[i].insertText(\"\u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the invoice document:\n#  - document date 27.5.2024 -> 6.6.2024 (both occurrences)\n#  - item name KALTMANN MASCHINEN 2 -> YH9000AE\n#  - remove the \"14450.00\" price/sum figures (all three occurrences)\n#  - shorten the spelled-out sum text to just \"\u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\"\n\n$d = $word.ActiveDocument\n\n# 1) Replace both occurrences of the old date with the new date.\nfor ($i = 0; $i -lt 2; $i++) {\n    $rngDate = $d.Content\n    $rngDate.Find.ClearFormatting()\n    if ($rngDate.Find.Execute(\"27.5.2024\")) {\n        $rngDate.Text = \"6.6.2024\"\n    }\n}\n\n# 2) Replace the equipment name.\n$rngName = $d.Content\n$rngName.Find.ClearFormatting()\nif ($rngName.Find.Execute(\"KALTMANN MASCHINEN 2\")) {\n    $rngName.Text = \"YH9000AE\"\n}\n\n# 3) Clear the \"14450.00\" amounts (unit price, line total, and the amount\n#    embedded in the \"\u0423\u0421\u042c\u041e\u0413\u041e\" summary cell).\nfor ($i = 0; $i -lt 3; $i++) {\n    $rngAmount = $d.Content\n    $rngAmount.Find.ClearFormatting()\n    if ($rngAmount.Find.Execute(\"14450.00\")) {\n        $rngAmount.Text = \"\"\n    }\n}\n\n# 4) Shorten the spelled-out sum text.\n$rngSum = $d.Content\n$rngSum.Find.ClearFormatting()\nif ($rngSum.Find.Execute(\"\u0447\u043e\u0442\u0438\u0440\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0447\u043e\u0442\u0438\u0440\u0438\u0441\u0442\u0430 \u043f\u044f\u0442\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\")) {\n    $rngSum.Text = \"\u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\"\n}\n"}
